$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data block ends at row 942 (last real "TESTE 2" row).
# Append 43 new test rows (943-985) following the same pattern as rows 941/942:
#   columns A, B, D, E, F, G, H = "TESTE"
#   column C = "TESTE 3" .. "TESTE 45"
$startRow = 943
$endTeste = 45

for ($n = 3; $n -le $endTeste; $n++) {
    $r = $startRow + ($n - 3)
    $ws.Cells.Item($r, 1).Value = "TESTE"
    $ws.Cells.Item($r, 2).Value = "TESTE"
    $ws.Cells.Item($r, 3).Value = "TESTE $n"
    $ws.Cells.Item($r, 4).Value = "TESTE"
    $ws.Cells.Item($r, 5).Value = "TESTE"
    $ws.Cells.Item($r, 6).Value = "TESTE"
    $ws.Cells.Item($r, 7).Value = "TESTE"
    $ws.Cells.Item($r, 8).Value = "TESTE"
}

$lastRow = $startRow + ($endTeste - 3)

# The source rows in columns E:G inherit column-level styling; clear it so the
# newly added cells stay unstyled like rows 941/942 did.
$ws.Range("E" + $startRow + ":G" + $lastRow).Style = "Normal"

# Grow the table (ListObject) to include the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J" + $lastRow))

# Update the view to match: scroll down a bit and move the selection to E954.
$win = $excel.Windows.Item(1)
$win.ScrollRow = 930
$ws.Range("E954").Select()
